$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (16) down to the new row (17)
$ws.Range("A16:T16").Copy()
$ws.Range("A17:T17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New local extreme: Kotlas
$ws.Range("A17").Value = 11710000
$ws.Range("B17").Value = "Котлас"
$ws.Range("C17").Value = 2019
$ws.Range("D17").Value = 0.29761904761904762
$ws.Range("E17").Value = 0.49603633085117038
$ws.Range("F17").Value = 0.4803921914504613
$ws.Range("G17").Value = 0.67446195458032387
$ws.Range("H17").Value = 0.81617276227703461
$ws.Range("I17").Value = 0.70416771217697205
$ws.Range("J17").Value = 0.56779397983662483
$ws.Range("K17").Value = 0.43217869860796371
$ws.Range("L17").Value = 0.18432592912931001
$ws.Range("M17").Value = 0.1690738406599411
$ws.Range("N17").Value = 0.029681448438395701
$ws.Range("O17").Value = 0.000068280893655534139
$ws.Range("P17").Value = 0.013062384580036
$ws.Range("Q17").Value = 0.0060770196999548002
$ws.Range("R17").Value = 0.063640035544066498
$ws.Range("S17").Value = 0.56576495966639262
$ws.Range("T17").Value = 0.094227873614833799

# Update the active selection to match the post-edit state
$ws.Range("D18").Select()
